$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "darsh2@gmail.com"
$ws.Range("A2").Value = "sanj2@gmail.com"
$ws.Range("A3").Value = "harshi2@gmail.com"

$ws.Range("A3").Select()
